$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 6, pushing existing rows 6-12
# down to 8-14 (matches the target dimension A1:T14).
$ws.Rows("6:7").Insert()

# New row 6 (price update for 2022-10-14, "Especial" quality)
$ws.Range("A6").Value = 7
$ws.Range("B6").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C6").Value = "Ñuble"
$ws.Range("D6").Value = 44848
$ws.Range("E6").Value = 16
$ws.Range("F6").Value = "Fruta"
$ws.Range("G6").Value = 100107
$ws.Range("H6").Value = "Otros"
$ws.Range("I6").Value = 100107002
$ws.Range("J6").Value = "Chirimoya"
$ws.Range("K6").Value = "Cultivar IV Región"
$ws.Range("L6").Value = "Especial"
$ws.Range("M6").Value = 60
$ws.Range("N6").Value = 24000
$ws.Range("O6").Value = 25000
$ws.Range("P6").Value = 24500
$ws.Range("Q6").Value = "`$/bandeja 10 kilos"
$ws.Range("R6").Value = "Provincia de Limarí"
$ws.Range("S6").Value = 2450
$ws.Range("T6").Value = 10

# New row 7 (price update for 2022-10-14, "Primera" quality)
$ws.Range("A7").Value = 7
$ws.Range("B7").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C7").Value = "Ñuble"
$ws.Range("D7").Value = 44848
$ws.Range("E7").Value = 16
$ws.Range("F7").Value = "Fruta"
$ws.Range("G7").Value = 100107
$ws.Range("H7").Value = "Otros"
$ws.Range("I7").Value = 100107002
$ws.Range("J7").Value = "Chirimoya"
$ws.Range("K7").Value = "Cultivar IV Región"
$ws.Range("L7").Value = "Primera"
$ws.Range("M7").Value = 120
$ws.Range("N7").Value = 21000
$ws.Range("O7").Value = 22000
$ws.Range("P7").Value = 21500
$ws.Range("Q7").Value = "`$/bandeja 10 kilos"
$ws.Range("R7").Value = "Provincia de Limarí"
$ws.Range("S7").Value = 2150
$ws.Range("T7").Value = 10
